$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card23")

# Add new header "Correction" in column N (14), row 1.
# Copy M1's format (bold font, border, centered alignment) onto N1 first so the
# new header cell shares the existing header style, then set its text.
$ws.Range("M1").Copy($ws.Range("N1"))
$ws.Range("N1").Value = "Correction"

# Column M (13) is currently blank (empty inline strings) for rows 2-12.
# Propagate that same "blank cell" shape into the new column N (14) first,
# before M gets populated with "nan" below.
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 13).Copy($ws.Cells.Item($r, 14))
}

# Now fill the previously-blank M2:M12 cells with "nan".
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 13).Value = "nan"
}
